{"js": "// Replace the date heading and every two-digit x two-digit multiplication\n// answer cell in the table with the new values from the updated worksheet.\nconst replacements = [\n  [\"2024-11-11 Monday\", \"2024-11-12 Tuesday\"],\n  [\"15\u00d751=765\", \"94\u00d791=8554\"],\n  [\"50\u00d795=4750\", \"67\u00d742=2814\"],\n  [\"76\u00d717=1292\", \"82\u00d767=5494\"],\n  [\"78\u00d727=2106\", \"89\u00d731=2759\"],\n  [\"32\u00d713=416\", \"50\u00d715=750\"],\n  [\"76\u00d715=1140\", \"64\u00d755=3520\"],\n  [\"42\u00d736=1512\", \"91\u00d721=1911\"],\n  [\"20\u00d780=1600\", \"41\u00d739=1599\"],\n  [\"65\u00d766=4290\", \"99\u00d751=5049\"],\n  [\"43\u00d769=2967\", \"91\u00d756=5096\"],\n  [\"33\u00d791=3003\", \"20\u00d763=1260\"],\n  [\"56\u00d757=3192\", \"27\u00d763=1701\"],\n  [\"28\u00d780=2240\", \"62\u00d782=5084\"],\n  [\"49\u00d763=3087\", \"96\u00d794=9024\"],\n  [\"51\u00d720=1020\", \"46\u00d716=736\"],\n  [\"81\u00d769=5589\", \"72\u00d729=2088\"],\n  [\"82\u00d756=4592\", \"46\u00d796=4416\"],\n  [\"18\u00d760=1080\", \"79\u00d761=4819\"],\n  [\"78\u00d796=7488\", \"92\u00d769=6348\"],\n  [\"64\u00d748=3072\", \"51\u00d740=2040\"],\n  [\"81\u00d759=4779\", \"49\u00d771=3479\"],\n  [\"22\u00d760=1320\", \"21\u00d740=840\"],\n  [\"72\u00d744=3168\", \"77\u00d743=3311\"],\n  [\"84\u00d724=2016\", \"18\u00d720=360\"],\n  [\"33\u00d794=3102\", \"47\u00d794=4418\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const hit of results.items) {\n    hit.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date heading and every two-digit x two-digit multiplication\n# answer cell in the table to the new values from the updated worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-11 Monday\", \"2024-11-12 Tuesday\"),\n    @(\"15\u00d751=765\", \"94\u00d791=8554\"),\n    @(\"50\u00d795=4750\", \"67\u00d742=2814\"),\n    @(\"76\u00d717=1292\", \"82\u00d767=5494\"),\n    @(\"78\u00d727=2106\", \"89\u00d731=2759\"),\n    @(\"32\u00d713=416\", \"50\u00d715=750\"),\n    @(\"76\u00d715=1140\", \"64\u00d755=3520\"),\n    @(\"42\u00d736=1512\", \"91\u00d721=1911\"),\n    @(\"20\u00d780=1600\", \"41\u00d739=1599\"),\n    @(\"65\u00d766=4290\", \"99\u00d751=5049\"),\n    @(\"43\u00d769=2967\", \"91\u00d756=5096\"),\n    @(\"33\u00d791=3003\", \"20\u00d763=1260\"),\n    @(\"56\u00d757=3192\", \"27\u00d763=1701\"),\n    @(\"28\u00d780=2240\", \"62\u00d782=5084\"),\n    @(\"49\u00d763=3087\", \"96\u00d794=9024\"),\n    @(\"51\u00d720=1020\", \"46\u00d716=736\"),\n    @(\"81\u00d769=5589\", \"72\u00d729=2088\"),\n    @(\"82\u00d756=4592\", \"46\u00d796=4416\"),\n    @(\"18\u00d760=1080\", \"79\u00d761=4819\"),\n    @(\"78\u00d796=7488\", \"92\u00d769=6348\"),\n    @(\"64\u00d748=3072\", \"51\u00d740=2040\"),\n    @(\"81\u00d759=4779\", \"49\u00d771=3479\"),\n    @(\"22\u00d760=1320\", \"21\u00d740=840\"),\n    @(\"72\u00d744=3168\", \"77\u00d743=3311\"),\n    @(\"84\u00d724=2016\", \"18\u00d720=360\"),\n    @(\"33\u00d794=3102\", \"47\u00d794=4418\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $true, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
